$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire row 681 (the "「イードおめでとう」" entry).
# This shifts all subsequent rows up by one, matching the target diff.
$ws.Rows.Item(681).EntireRow.Delete()
